# Auto-generated edit script applying numeric corrections to Sheets (Hades_Profits workbook).
# For each affected row in ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR, update H:N profit-calculation cells
# to the refreshed values produced by the scheduled runner, deleting any cells that no longer apply.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 56.5
$ws.Range("I2").Value = 30.4
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 30.4
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 82.59999999999999
$ws.Range("N2").Value = -326
$ws.Range("H9").Value = 264.2857
$ws.Range("I9").Value = 180
$ws.Range("J9").Value = 475
$ws.Range("K9").Value = 180
$ws.Range("L9").Value = 475
$ws.Range("M9").Value = -11
$ws.Range("N9").Value = -813
$ws.Range("H32").Value = 817.4545000000001
$ws.Range("I32").Value = 832.6667
$ws.Range("J32").Value = 811.75
$ws.Range("K32").Value = 832.6667
$ws.Range("L32").Value = 811.75
$ws.Range("M32").Value = -506.6667
$ws.Range("N32").Value = -1463.75
$ws.Range("H98").Value = 1185
$ws.Range("I98").Value = 1182
$ws.Range("J98").Value = 1200
$ws.Range("K98").Value = 1182
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = 316
$ws.Range("N98").Value = -4196
$ws.Range("H113").Value = 3072.7666
$ws.Range("I113").Value = 2694.682
$ws.Range("K113").Value = 2694.682
$ws.Range("M113").Value = 559.3180000000002
$ws.Range("H116").Value = 1940.2632
$ws.Range("I116").Value = 1776.0714
$ws.Range("K116").Value = 1776.0714
$ws.Range("M116").Value = 1665.9286
$ws.Range("H122").Value = 1185
$ws.Range("I122").Value = 1182
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 3546
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1096
$ws.Range("N122").Value = -8500
$ws.Range("H127").Value = 1659.5834
$ws.Range("I127").Value = 599.3333
$ws.Range("J127").Value = 2013
$ws.Range("K127").Value = 1797.9999
$ws.Range("L127").Value = 6039
$ws.Range("M127").Value = 3162.0001
$ws.Range("N127").Value = -15959
$ws.Range("H132").Value = 861460.7
$ws.Range("I132").Value = 1634.2979
$ws.Range("K132").Value = 4902.893700000001
$ws.Range("M132").Value = -2372.893700000001
$ws.Range("H135").Value = 22040.041
$ws.Range("I135").Value = 28230.973
$ws.Range("J135").Value = 2951.3333
$ws.Range("K135").Value = 254078.757
$ws.Range("L135").Value = 26561.9997
$ws.Range("M135").Value = -251543.757
$ws.Range("N135").Value = -31631.9997
$ws.Range("H136").Value = 37113.332
$ws.Range("J136").Value = 37113.332
$ws.Range("L136").Value = 37113.332
$ws.Range("N136").Value = -47313.332
$ws.Range("H137").Value = 2633303
$ws.Range("I137").Value = 5001518
$ws.Range("J137").Value = 1952.7778
$ws.Range("K137").Value = 15004554
$ws.Range("L137").Value = 5858.3334
$ws.Range("M137").Value = -15002004
$ws.Range("N137").Value = -10958.3334
$ws.Range("H141").Value = 520.1053000000001
$ws.Range("I141").Value = 520.1053000000001
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1560.3159
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3619.6841
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 20450000
$ws.Range("I61").Value = 22751000
$ws.Range("K61").Value = 22751000
$ws.Range("M61").Value = -22750788
$ws.Range("H74").Value = 6633107.5
$ws.Range("I74").Value = 8965603
$ws.Range("J74").Value = 102120
$ws.Range("K74").Value = 8965603
$ws.Range("L74").Value = 102120
$ws.Range("M74").Value = -8964729
$ws.Range("N74").Value = -103868
$ws.Range("H77").Value = 6633107.5
$ws.Range("I77").Value = 8965603
$ws.Range("J77").Value = 102120
$ws.Range("K77").Value = 44828015
$ws.Range("L77").Value = 510600
$ws.Range("M77").Value = -44823647
$ws.Range("N77").Value = -519336
$ws.Range("H122").Value = 4832720.5
$ws.Range("I122").Value = 2072.1333
$ws.Range("J122").Value = 13890186
$ws.Range("K122").Value = 6216.3999
$ws.Range("L122").Value = 41670558
$ws.Range("M122").Value = -3766.3999
$ws.Range("N122").Value = -41675458
$ws.Range("H132").Value = 44327.457
$ws.Range("I132").Value = 27820.541
$ws.Range("J132").Value = 99850.73
$ws.Range("K132").Value = 83461.62300000001
$ws.Range("L132").Value = 299552.19
$ws.Range("M132").Value = -80931.62300000001
$ws.Range("N132").Value = -304612.19
$ws.Range("H136").Value = 20450000
$ws.Range("I136").Value = 22751000
$ws.Range("K136").Value = 68253000
$ws.Range("M136").Value = -68250450
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1508.4509
$ws.Range("I134").Value = 884.9773
$ws.Range("J134").Value = 5427.4287
$ws.Range("K134").Value = 2654.9319
$ws.Range("L134").Value = 16282.2861
$ws.Range("M134").Value = -119.9319
$ws.Range("N134").Value = -21352.2861
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 83333840
$ws.Range("J22").Value = 714.2857
$ws.Range("L22").Value = 714.2857
$ws.Range("N22").Value = -1414.2857
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5480
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("N27").Value = -5384
$ws.Range("H31").Value = 3237.4324
$ws.Range("I31").Value = 1185.6897
$ws.Range("J31").Value = 10675
$ws.Range("K31").Value = 1185.6897
$ws.Range("L31").Value = 10675
$ws.Range("M31").Value = -890.6896999999999
$ws.Range("N31").Value = -11265
$ws.Range("H34").Value = 3237.4324
$ws.Range("I34").Value = 1185.6897
$ws.Range("J34").Value = 10675
$ws.Range("K34").Value = 1185.6897
$ws.Range("L34").Value = 10675
$ws.Range("M34").Value = -983.6896999999999
$ws.Range("N34").Value = -11079
$ws.Range("H50").Value = 27425.334
$ws.Range("J50").Value = 27425.334
$ws.Range("L50").Value = 27425.334
$ws.Range("N50").Value = -28675.334
$ws.Range("H51").Value = 28099
$ws.Range("J51").Value = 28099
$ws.Range("L51").Value = 28099
$ws.Range("N51").Value = -29571
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 28099
$ws.Range("J61").Value = 28099
$ws.Range("L61").Value = 28099
$ws.Range("N61").Value = -28795
$ws.Range("H99").Value = 10333
$ws.Range("I99").Value = 9999
$ws.Range("K99").Value = 9999
$ws.Range("M99").Value = -8501
$ws.Range("H122").Value = 1776.0454
$ws.Range("I122").Value = 1559.9375
$ws.Range("J122").Value = 2352.3333
$ws.Range("K122").Value = 4679.8125
$ws.Range("L122").Value = 7056.999899999999
$ws.Range("M122").Value = -2229.8125
$ws.Range("N122").Value = -11956.9999
$ws.Range("H126").Value = 10333
$ws.Range("I126").Value = 9999
$ws.Range("K126").Value = 29997
$ws.Range("M126").Value = -27527
$ws.Range("H132").Value = 40014.617
$ws.Range("I132").Value = 26567
$ws.Range("J132").Value = 84840
$ws.Range("K132").Value = 79701
$ws.Range("L132").Value = 254520
$ws.Range("M132").Value = -77171
$ws.Range("N132").Value = -259580
$ws.Range("H134").Value = 32407.4
$ws.Range("I134").Value = 1609.625
$ws.Range("K134").Value = 4828.875
$ws.Range("M134").Value = -2293.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2000
$ws.Range("J113").Value = 2000
$ws.Range("L113").Value = 6000
$ws.Range("N113").Value = -10340
$ws.Range("H131").Value = 850.8
$ws.Range("I131").Value = 664.75
$ws.Range("J131").Value = 864.0893
$ws.Range("K131").Value = 1994.25
$ws.Range("L131").Value = 2592.2679
$ws.Range("M131").Value = 3045.75
$ws.Range("N131").Value = -12672.2679
$ws.Range("H132").Value = 2400.4546
$ws.Range("I132").Value = 4450
$ws.Range("K132").Value = 40050
$ws.Range("M132").Value = -37520
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 52683.18
$ws.Range("I132").Value = 31678.121
$ws.Range("K132").Value = 95034.363
$ws.Range("M132").Value = -92504.363
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H132").Value = 31517.318
$ws.Range("I132").Value = 20685.186
$ws.Range("J132").Value = 70513
$ws.Range("K132").Value = 62055.558
$ws.Range("L132").Value = 211539
$ws.Range("M132").Value = -59525.558
$ws.Range("N132").Value = -216599
$ws.Range("H136").Value = 83635.266
$ws.Range("I136").Value = 74211.47
$ws.Range("J136").Value = 96485.91
$ws.Range("K136").Value = 222634.41
$ws.Range("L136").Value = 289457.73
$ws.Range("M136").Value = -220084.41
$ws.Range("N136").Value = -294557.73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 37041.695
$ws.Range("I132").Value = 29125.844
$ws.Range("K132").Value = 87377.53200000001
$ws.Range("M132").Value = -84847.53200000001
